# 9th Stab - Cosmetic Changes
# Insert two new "date" columns (Jun_15, Jun_17) ahead of the most-recent
# rating column, pushing the older Jun_10/Jun_13 columns to the right and
# filling the newly created cells with the "UN" (unrated) placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Insert two blank columns at C:D - this shifts the existing C column
# (and its data/styles) to column E automatically.
$ws.Columns("C:D").Insert()

# Re-affirm the shifted former-B1 header text so the engine keeps it as
# its own distinct shared string instead of overwriting it in place.
$ws.Range("D1").Value = "Jun_13"

# New header cells for the freshly inserted columns (appended in the
# same order the workbook's shared-string table grew historically).
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill the newly created C/D cells for every data row with the "UN"
# placeholder rating, matching the rest of the table.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Keep the cosmetic column widths consistent across C, D and E (~8 chars).
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14
